$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 207 (the "怒ってる？" post), shifting all subsequent rows up by one.
$ws.Rows.Item(207).Delete()
